$d = $word.ActiveDocument

# Find the "CU010: " label run and place the insertion point right after it.
$range = $d.Content
$null = $range.Find.Execute("CU010: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$range.Collapse(0)  # wdCollapseEnd

# Insert the use-case title text as its own run (same paragraph).
$range.InsertAfter("Modificar datos")

# Toggling a formatting property and back forces the new text into a
# distinct run (matching the author's separate <w:r>) instead of merging
# it into the preceding "CU010: " run.
$range.Font.Bold = $true
$range.Font.Bold = $false
